$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.664.89'
$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.672.35'
$ws.Range("E3").Value = '  -1.96%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.99'
$ws.Range("E5").Value = '  -2.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.60'
$ws.Range("E6").Value = '  -2.05%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("E9").Value = '  -4.07%  '

$ws.Range("E10").Value = '  -3.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  -2.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  -4.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.147.73'
$ws.Range("E13").Value = '  -1.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.00'
$ws.Range("E14").Value = '  -2.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.600.50'
$ws.Range("E15").Value = '  -1.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000145'
$ws.Range("E16").Value = '  -3.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.678.95'
$ws.Range("E17").Value = '  -1.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.77'
$ws.Range("E18").Value = '  -6.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.57'
$ws.Range("E19").Value = '  -3.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.84'
$ws.Range("E20").Value = '  -3.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.17'
$ws.Range("E21").Value = '  -5.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.506'
$ws.Range("E23").Value = '  -2.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.10'
$ws.Range("E24").Value = '  -2.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.10'
$ws.Range("E27").Value = '  -3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.39'
$ws.Range("E28").Value = '  +3.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0844'
$ws.Range("E29").Value = '  -6.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.21'
$ws.Range("E30").Value = '  +0.57%  '

$ws.Range("E31").Value = '  -2.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.15'
$ws.Range("E32").Value = '  -2.50%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  -2.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.44'
$ws.Range("E35").Value = '  -2.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.36'
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.77'
$ws.Range("E37").Value = '  -1.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '335.62'
$ws.Range("E38").Value = '  -3.03%  '

$ws.Range("E39").Value = '  -2.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.917'
$ws.Range("E40").Value = '  -5.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.96'
$ws.Range("E41").Value = '  -3.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.25'
$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.72'
$ws.Range("E43").Value = '  -5.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.04'
$ws.Range("E44").Value = '  -4.71%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.613'
$ws.Range("E46").Value = '  -1.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0552'
$ws.Range("E47").Value = '  -5.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.01'
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.34'
$ws.Range("E49").Value = '  -1.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0965'
$ws.Range("E50").Value = '  -3.33%  '

$ws.Range("E51").Value = '  -4.55%  '
